$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'69.718.50"
$ws.Range("E2").Value = "  -0.29%  "

# Row 3
$ws.Range("D3").Value = "'3.526.90"
$ws.Range("E3").Value = "  +0.59%  "

# Row 4
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  +0.07%  "

# Row 5
$ws.Range("D5").Value = "'606.34"
$ws.Range("E5").Value = "  -0.02%  "

# Row 6
$ws.Range("D6").Value = "'194.92"
$ws.Range("E6").Value = "  +1.10%  "

# Row 7
$ws.Range("D7").Value = "'0.624"
$ws.Range("E7").Value = "  -0.60%  "

# Row 9
$ws.Range("D9").Value = "'0.202"
$ws.Range("E9").Value = "  -6.03%  "

# Row 10
$ws.Range("E10").Value = "  -2.55%  "

# Row 11
$ws.Range("D11").Value = "'53.38"
$ws.Range("E11").Value = "  -0.41%  "

# Row 12
$ws.Range("D12").Value = "'0.0000301"
$ws.Range("E12").Value = "  -2.19%  "

# Row 13
$ws.Range("E13").Value = "  -1.76%  "

# Row 14
$ws.Range("D14").Value = "'4.090.52"
$ws.Range("E14").Value = "  +0.54%  "

# Row 15
$ws.Range("D15").Value = "'592.58"
$ws.Range("E15").Value = "  -4.02%  "

# Row 16
$ws.Range("D16").Value = "'69.864.91"
$ws.Range("E16").Value = "  -0.19%  "

# Row 17
$ws.Range("D17").Value = "'12.69"
$ws.Range("E17").Value = "  -0.14%  "

# Row 18
$ws.Range("D18").Value = "'18.96"
$ws.Range("E18").Value = "  +0.34%  "

# Row 19
$ws.Range("D19").Value = "'3.522.35"
$ws.Range("E19").Value = "  +0.38%  "

# Row 20
$ws.Range("E20").Value = "  +1.83%  "

# Row 21
$ws.Range("D21").Value = "'0.982"
$ws.Range("E21").Value = "  -1.12%  "

# Row 22
$ws.Range("D22").Value = "'17.78"
$ws.Range("E22").Value = "  -0.97%  "

# Row 23
$ws.Range("D23").Value = "'5.13"
$ws.Range("E23").Value = "  +2.21%  "

# Row 24
$ws.Range("D24").Value = "'102.61"
$ws.Range("E24").Value = "  -2.61%  "

# Row 25
$ws.Range("D25").Value = "'4.64"
$ws.Range("E25").Value = "  -0.21%  "

# Row 26
$ws.Range("D26").Value = "'3.04"
$ws.Range("E26").Value = "  -0.13%  "

# Row 27
$ws.Range("D27").Value = "'10.76"
$ws.Range("E27").Value = "  -2.26%  "

# Row 28
$ws.Range("D28").Value = "'9.51"
$ws.Range("E28").Value = "  -3.88%  "

# Row 29
$ws.Range("D29").Value = "'33.13"
$ws.Range("E29").Value = "  -3.74%  "

# Row 30
$ws.Range("D30").Value = "'7.04"
$ws.Range("E30").Value = "  -1.44%  "

# Row 31
$ws.Range("D31").Value = "'4.22"
$ws.Range("E31").Value = "  -1.44%  "

# Row 32
$ws.Range("D32").Value = "'12.32"
$ws.Range("E32").Value = "  -2.65%  "

# Row 33
$ws.Range("E33").Value = "  -0.44%  "

# Row 34
$ws.Range("D34").Value = "'63.37"
$ws.Range("E34").Value = "  -1.07%  "

# Row 35
$ws.Range("B35").Value = "Maker"
$ws.Range("C35").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D35").Value = "'3.796.49"
$ws.Range("E35").Value = "  +1.51%  "

# Row 36
$ws.Range("B36").Value = "Fetch.AI"
$ws.Range("C36").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D36").Value = "'3.17"
$ws.Range("E36").Value = "  +2.47%  "

# Row 37
$ws.Range("B37").Value = "Dai"
$ws.Range("C37").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D37").Value = "'1.00"
$ws.Range("E37").Value = "  +0.18%  "

# Row 38
$ws.Range("B38").Value = "PEPE"
$ws.Range("C38").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D38").Value = "'0.0₃0805"
$ws.Range("E38").Value = "  +1.07%  "

# Row 39
$ws.Range("D39").Value = "'514.27"
$ws.Range("E39").Value = "  -2.41%  "

# Row 40
$ws.Range("E40").Value = "  -0.12%  "

# Row 41
$ws.Range("E41").Value = "  -0.28%  "

# Row 42
$ws.Range("D42").Value = "'36.47"
$ws.Range("E42").Value = "  -1.00%  "

# Row 43
$ws.Range("E43").Value = "  -3.06%  "

# Row 44
$ws.Range("D44").Value = "'0.0447"
$ws.Range("E44").Value = "  -3.54%  "

# Row 45
$ws.Range("D45").Value = "'0.139"
$ws.Range("E45").Value = "  -1.33%  "

# Row 46
$ws.Range("D46").Value = "'2.81"
$ws.Range("E46").Value = "  -2.33%  "

# Row 47
$ws.Range("D47").Value = "'3.28"
$ws.Range("E47").Value = "  -1.23%  "

# Row 48
$ws.Range("E48").Value = "  +0.10%  "

# Row 49
$ws.Range("D49").Value = "'8.46"
$ws.Range("E49").Value = "  -3.40%  "

# Row 50
$ws.Range("B50").Value = "FLOKI"
$ws.Range("C50").Value = "https://coinranking.com/coin/fmHk13Rqw+floki-floki"
$ws.Range("D50").Value = "'0.000247"
$ws.Range("E50").Value = "  +3.80%  "

# Row 51
$ws.Range("B51").Value = "Mantle"
$ws.Range("C51").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D51").Value = "'1.32"
$ws.Range("E51").Value = "  +1.32%  "
